$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 605.6667
$ws.Range("I4").Value = 168.66667
$ws.Range("K4").Value = 168.66667
$ws.Range("M4").Value = -54.66667000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 545.3
$ws.Range("I33").Value = 424
$ws.Range("J33").Value = 1232.6666
$ws.Range("K33").Value = 424
$ws.Range("L33").Value = 1232.6666
$ws.Range("M33").Value = -195
$ws.Range("N33").Value = -1690.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11877.6
$ws.Range("I51").Value = 6999
$ws.Range("J51").Value = 13097.25
$ws.Range("K51").Value = 6999
$ws.Range("L51").Value = 13097.25
$ws.Range("M51").Value = -6515
$ws.Range("N51").Value = -14065.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 27785022
$ws.Range("I64").Value = 47622896
$ws.Range("K64").Value = 47622896
$ws.Range("M64").Value = -47622648

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 27785022
$ws.Range("I67").Value = 47622896
$ws.Range("K67").Value = 47622896
$ws.Range("M67").Value = -47622038

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 8379.6
$ws.Range("J96").Value = 11299.667
$ws.Range("L96").Value = 33899.001
$ws.Range("N96").Value = -36645.001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 3389.7
$ws.Range("I99").Value = 438.14285
$ws.Range("J99").Value = 10276.667
$ws.Range("K99").Value = 1314.42855
$ws.Range("L99").Value = 30830.001
$ws.Range("M99").Value = 183.5714499999999
$ws.Range("N99").Value = -33826.001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 897.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 897.5
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").Value = 2692.5
$ws.Range("N103").Value = -3864.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3602.25
$ws.Range("J112").Value = 2937.9412
$ws.Range("L112").Value = 8813.8236
$ws.Range("N112").Value = -11029.8236

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3996.3333
$ws.Range("I137").Value = 3878.9443
$ws.Range("J137").Value = 4172.4165
$ws.Range("K137").Value = 11636.8329
$ws.Range("L137").Value = 12517.2495
$ws.Range("M137").Value = -9086.832900000001
$ws.Range("N137").Value = -17617.2495

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3220
$ws.Range("I138").Value = 1388.3529
$ws.Range("J138").Value = 3855.4695
$ws.Range("K138").Value = 4165.0587
$ws.Range("L138").Value = 11566.4085
$ws.Range("M138").Value = 974.9412999999995
$ws.Range("N138").Value = -21846.4085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3092.3872
$ws.Range("I61").Value = 2279.2808
$ws.Range("J61").Value = 4379.8057
$ws.Range("K61").Value = 2279.2808
$ws.Range("L61").Value = 4379.8057
$ws.Range("M61").Value = -2067.2808
$ws.Range("N61").Value = -4803.8057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4492.8335
$ws.Range("I74").Value = 4929.4375
$ws.Range("K74").Value = 4929.4375
$ws.Range("M74").Value = -4055.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 48809.445
$ws.Range("J76").Value = 48809.445
$ws.Range("L76").Value = 48809.445
$ws.Range("N76").Value = -49485.445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4492.8335
$ws.Range("I77").Value = 4929.4375
$ws.Range("K77").Value = 24647.1875
$ws.Range("M77").Value = -20279.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 48809.445
$ws.Range("J79").Value = 48809.445
$ws.Range("L79").Value = 48809.445
$ws.Range("N79").Value = -51149.445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5389.9414
$ws.Range("I122").Value = 4013.3333
$ws.Range("K122").Value = 12039.9999
$ws.Range("M122").Value = -9589.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 710970.8
$ws.Range("I132").Value = 837032.9399999999
$ws.Range("J132").Value = 105872.6
$ws.Range("K132").Value = 2511098.82
$ws.Range("L132").Value = 317617.8
$ws.Range("M132").Value = -2508568.82
$ws.Range("N132").Value = -322677.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3092.3872
$ws.Range("I136").Value = 2279.2808
$ws.Range("J136").Value = 4379.8057
$ws.Range("K136").Value = 6837.8424
$ws.Range("L136").Value = 13139.4171
$ws.Range("M136").Value = -4287.8424
$ws.Range("N136").Value = -18239.4171

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4744.4517
$ws.Range("I94").Value = 2547.55
$ws.Range("J94").Value = 8738.817999999999
$ws.Range("K94").Value = 2547.55
$ws.Range("L94").Value = 8738.817999999999
$ws.Range("M94").Value = -2096.55
$ws.Range("N94").Value = -9640.817999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 44975
$ws.Range("J95").Value = 44975
$ws.Range("L95").Value = 44975
$ws.Range("N95").Value = -50467

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7913.7607
$ws.Range("I31").Value = 20298.2
$ws.Range("J31").Value = 4473.6387
$ws.Range("K31").Value = 20298.2
$ws.Range("L31").Value = 4473.6387
$ws.Range("M31").Value = -20003.2
$ws.Range("N31").Value = -5063.6387

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7913.7607
$ws.Range("I34").Value = 20298.2
$ws.Range("J34").Value = 4473.6387
$ws.Range("K34").Value = 20298.2
$ws.Range("L34").Value = 4473.6387
$ws.Range("M34").Value = -20096.2
$ws.Range("N34").Value = -4877.6387

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 69270
$ws.Range("J43").Value = 69270
$ws.Range("L43").Value = 69270
$ws.Range("N43").Value = -69638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8405.799999999999
$ws.Range("I62").Value = 8149.2
$ws.Range("K62").Value = 8149.2
$ws.Range("M62").Value = -7525.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 8405.799999999999
$ws.Range("I65").Value = 8149.2
$ws.Range("K65").Value = 40746
$ws.Range("M65").Value = -37626

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 69270
$ws.Range("J101").Value = 69270
$ws.Range("L101").Value = 69270
$ws.Range("N101").Value = -75760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 7292.143
$ws.Range("I122").Value = 3107.0908
$ws.Range("K122").Value = 9321.2724
$ws.Range("M122").Value = -6871.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7702.6895
$ws.Range("I132").Value = 3559.9565
$ws.Range("K132").Value = 10679.8695
$ws.Range("M132").Value = -8149.869499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 82.57143000000001
$ws.Range("I2").Value = 50
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 27856.715
$ws.Range("I87").Value = 24999
$ws.Range("K87").Value = 74997
$ws.Range("M87").Value = -73749

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 27856.715
$ws.Range("I90").Value = 24999
$ws.Range("K90").Value = 224991
$ws.Range("M90").Value = -218751

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 200
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 200
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("M92").Value = 600
$ws.Range("N92").Value = -3096

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 28500
$ws.Range("J106").Value = 28500
$ws.Range("L106").Value = 85500
$ws.Range("N106").Value = -87392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 22782.555
$ws.Range("I110").Value = 17008.6
$ws.Range("K110").Value = 51025.8
$ws.Range("M110").Value = -46935.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2841852.2
$ws.Range("I113").Value = 7812908
$ws.Range("J113").Value = 1248.9286
$ws.Range("K113").Value = 23438724
$ws.Range("L113").Value = 3746.7858
$ws.Range("M113").Value = -23436554
$ws.Range("N113").Value = -8086.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3009.7727
$ws.Range("J122").Value = 3299.2632
$ws.Range("L122").Value = 29693.3688
$ws.Range("N122").Value = -34593.3688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2222.647
$ws.Range("I129").Value = 472.66666
$ws.Range("J129").Value = 4191.375
$ws.Range("K129").Value = 1417.99998
$ws.Range("L129").Value = 12574.125
$ws.Range("M129").Value = 3582.00002
$ws.Range("N129").Value = -22574.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 17821.5
$ws.Range("J98").Value = 17821.5
$ws.Range("L98").Value = 17821.5
$ws.Range("N98").Value = -23811.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7281.067
$ws.Range("I122").Value = 5594.8696
$ws.Range("K122").Value = 16784.6088
$ws.Range("M122").Value = -14334.6088

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5049.9556
$ws.Range("I132").Value = 5016.282
$ws.Range("K132").Value = 15048.846
$ws.Range("M132").Value = -12518.846

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8539.933999999999
$ws.Range("I7").Value = 5623.231
$ws.Range("J7").Value = 27498.5
$ws.Range("K7").Value = 5623.231
$ws.Range("L7").Value = 27498.5
$ws.Range("M7").Value = -5511.231
$ws.Range("N7").Value = -27722.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 33335578
$ws.Range("I46").Value = 1200
$ws.Range("K46").Value = 1200
$ws.Range("M46").Value = -1012

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 18754
$ws.Range("J101").Value = 18754
$ws.Range("L101").Value = 18754
$ws.Range("N101").Value = -25244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 73495
$ws.Range("J105").Value = 73495
$ws.Range("L105").Value = 73495
$ws.Range("N105").Value = -80483

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8539.933999999999
$ws.Range("I126").Value = 5623.231
$ws.Range("J126").Value = 27498.5
$ws.Range("K126").Value = 16869.693
$ws.Range("L126").Value = 82495.5
$ws.Range("M126").Value = -14399.693
$ws.Range("N126").Value = -87435.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3368.9
$ws.Range("J132").Value = 1665
$ws.Range("L132").Value = 4995
$ws.Range("N132").Value = -10055

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5463.763
$ws.Range("I136").Value = 5084.0713
$ws.Range("J136").Value = 6526.9
$ws.Range("K136").Value = 15252.2139
$ws.Range("L136").Value = 19580.7
$ws.Range("M136").Value = -12702.2139
$ws.Range("N136").Value = -24680.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 9403.333000000001
$ws.Range("I58").Value = 6500
$ws.Range("J58").Value = 15210
$ws.Range("K58").Value = 6500
$ws.Range("L58").Value = 15210
$ws.Range("M58").Value = -6192
$ws.Range("N58").Value = -15826
